# Apply updated Price (D) / Volume(1h) (E) figures from the latest cryptos data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    # Force the literal text into the cell (leading apostrophe keeps Excel from
    # auto-coercing decimal-looking strings to numbers), then strip the resulting
    # quote-prefix style back off so formatting matches the untouched cells.
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.856.64"
$ws.Range("E2").Value = "  +1.84%  "
Set-TextValue $ws.Range("D3") "3.459.60"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws.Range("D5") "576.07"
$ws.Range("E5").Value = "  -0.32%  "
Set-TextValue $ws.Range("D6") "160.44"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("E7").Value = "  +0.10%  "
Set-TextValue $ws.Range("D8") "3.461.19"
$ws.Range("E8").Value = "  +0.34%  "
Set-TextValue $ws.Range("D9") "0.580"
$ws.Range("E9").Value = "  +8.66%  "
Set-TextValue $ws.Range("D10") "7.36"
$ws.Range("E10").Value = "  -2.73%  "
$ws.Range("E11").Value = "  +2.30%  "
Set-TextValue $ws.Range("D12") "0.440"
$ws.Range("E12").Value = "  +0.67%  "
Set-TextValue $ws.Range("D13") "4.059.99"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("E14").Value = "  -2.50%  "
Set-TextValue $ws.Range("D15") "0.0000194"
$ws.Range("E15").Value = "  +4.06%  "
$ws.Range("E16").Value = "  +3.17%  "
Set-TextValue $ws.Range("D17") "64.901.57"
$ws.Range("E17").Value = "  +1.70%  "
Set-TextValue $ws.Range("D18") "3.470.68"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("E19").Value = "  -1.03%  "
Set-TextValue $ws.Range("D20") "14.30"
$ws.Range("E20").Value = "  +0.35%  "
Set-TextValue $ws.Range("D21") "388.89"
$ws.Range("E21").Value = "  -0.81%  "
Set-TextValue $ws.Range("D22") "8.20"
$ws.Range("E22").Value = "  -3.54%  "
Set-TextValue $ws.Range("D23") "73.17"
$ws.Range("E23").Value = "  +1.76%  "
Set-TextValue $ws.Range("D24") "0.544"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +14.59%  "
Set-TextValue $ws.Range("D27") "9.67"
$ws.Range("E27").Value = "  +1.33%  "
Set-TextValue $ws.Range("D28") "0.181"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("E29").Value = "  +0.01%  "
Set-TextValue $ws.Range("D30") "6.23"
$ws.Range("E30").Value = "  +8.32%  "
$ws.Range("E31").Value = "  +4.93%  "
$ws.Range("E32").Value = "  +0.08%  "
Set-TextValue $ws.Range("D33") "23.69"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +4.65%  "
$ws.Range("E39").Value = "  +2.65%  "
Set-TextValue $ws.Range("D40") "3.006.27"
$ws.Range("E40").Value = "  +2.51%  "
Set-TextValue $ws.Range("D41") "0.0766"
$ws.Range("E41").Value = "  -1.91%  "
Set-TextValue $ws.Range("D42") "27.26"
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("E43").Value = "  +4.86%  "
Set-TextValue $ws.Range("D44") "42.98"
$ws.Range("E44").Value = "  +2.83%  "
$ws.Range("E45").Value = "  -1.82%  "
Set-TextValue $ws.Range("D46") "0.775"
$ws.Range("E46").Value = "  +0.62%  "
Set-TextValue $ws.Range("D47") "24.38"
$ws.Range("E47").Value = "  +8.01%  "
$ws.Range("E48").Value = "  +0.53%  "
Set-TextValue $ws.Range("D49") "0.879"
$ws.Range("E49").Value = "  +7.14%  "
Set-TextValue $ws.Range("D50") "6.59"
$ws.Range("E50").Value = "  +3.23%  "
Set-TextValue $ws.Range("D51") "305.15"
$ws.Range("E51").Value = "  +2.90%  "

# Rows 37/38: ImmutableX and Monero swap rank positions, with refreshed figures.
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "1.51"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D38") "163.14"
$ws.Range("E38").Value = "  +2.81%  "
